$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (Recursos) updates -------------------------------------------------
# Rows whose resource text changes (existing "Enc Proyecto"/"Analista" entries are
# being replaced by actual people's names) and rows that get a brand-new resource.
$ws.Range("E2").Value = "Paolo"
$ws.Range("E3").Value = "Paolo"
$ws.Range("E4").Value = "Yuliana"
$ws.Range("E5").Value = "Paolo"
$ws.Range("E7").Value = " Paolo"
$ws.Range("E8").Value = "Paolo"
$ws.Range("E9").Value = "Paolo"
$ws.Range("E10").Value = "Diego"
$ws.Range("E12").Value = "Diego"
$ws.Range("E13").Value = "Tatiana"
$ws.Range("E14").Value = "Yuliana"
$ws.Range("E16").Value = "Yuliana"
$ws.Range("E18").Value = "Diego/Tatiana"
$ws.Range("E20").Value = "Yanela"
$ws.Range("E21").Value = "Yanela"
$ws.Range("E22").Value = "Diego"
$ws.Range("E23").Value = "Tatiana"
$ws.Range("E24").Value = "Yanela"

# --- Column C (Duración días) updates ---------------------------------------------
$ws.Range("C9").Value = 7
$ws.Range("C16").Value = 25
$ws.Range("C18").Value = 30
$ws.Range("C20").Value = 5
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 5
$ws.Range("C23").Value = 4

# --- New hidden helper column F ----------------------------------------------------
$ws.Columns("F").ColumnWidth = 0
$ws.Columns("F").Hidden = $true

# --- Selection / active cell --------------------------------------------------------
$ws.Range("C5").Select() | Out-Null
